{"js": "// Replace each two-digit multiplication equation in the table with its new value.\n// Each old equation string is unique in the document, so a simple search+replace\n// per pair is sufficient and avoids any row/index bookkeeping.\nconst replacements = [\n  ['11\u00d712=132', '48\u00d725=1200'],\n  ['59\u00d798=5782', '20\u00d784=1680'],\n  ['56\u00d734=1904', '27\u00d733=891'],\n  ['99\u00d785=8415', '17\u00d754=918'],\n  ['92\u00d744=4048', '44\u00d739=1716'],\n  ['53\u00d778=4134', '81\u00d796=7776'],\n  ['65\u00d768=4420', '34\u00d735=1190'],\n  ['66\u00d786=5676', '85\u00d791=7735'],\n  ['17\u00d747=799', '71\u00d758=4118'],\n  ['38\u00d721=798', '87\u00d756=4872'],\n  ['38\u00d712=456', '26\u00d750=1300'],\n  ['31\u00d764=1984', '68\u00d775=5100'],\n  ['30\u00d791=2730', '72\u00d724=1728'],\n  ['84\u00d732=2688', '27\u00d788=2376'],\n  ['28\u00d711=308', '68\u00d713=884'],\n  ['16\u00d738=608', '40\u00d732=1280'],\n  ['18\u00d790=1620', '25\u00d721=525'],\n  ['25\u00d716=400', '43\u00d734=1462'],\n  ['85\u00d740=3400', '30\u00d724=720'],\n  ['37\u00d711=407', '85\u00d735=2975'],\n  ['34\u00d746=1564', '15\u00d711=165'],\n  ['35\u00d747=1645', '94\u00d776=7144'],\n  ['42\u00d797=4074', '64\u00d736=2304'],\n  ['11\u00d734=374', '93\u00d771=6603'],\n  ['62\u00d774=4588', '27\u00d794=2538'],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation in the table with its new value.\n# Each old equation string is unique in the document, so Find/Replace (wdReplaceAll)\n# per pair is sufficient and avoids row/cell index bookkeeping.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"11\u00d712=132\", \"48\u00d725=1200\"),\n    @(\"59\u00d798=5782\", \"20\u00d784=1680\"),\n    @(\"56\u00d734=1904\", \"27\u00d733=891\"),\n    @(\"99\u00d785=8415\", \"17\u00d754=918\"),\n    @(\"92\u00d744=4048\", \"44\u00d739=1716\"),\n    @(\"53\u00d778=4134\", \"81\u00d796=7776\"),\n    @(\"65\u00d768=4420\", \"34\u00d735=1190\"),\n    @(\"66\u00d786=5676\", \"85\u00d791=7735\"),\n    @(\"17\u00d747=799\", \"71\u00d758=4118\"),\n    @(\"38\u00d721=798\", \"87\u00d756=4872\"),\n    @(\"38\u00d712=456\", \"26\u00d750=1300\"),\n    @(\"31\u00d764=1984\", \"68\u00d775=5100\"),\n    @(\"30\u00d791=2730\", \"72\u00d724=1728\"),\n    @(\"84\u00d732=2688\", \"27\u00d788=2376\"),\n    @(\"28\u00d711=308\", \"68\u00d713=884\"),\n    @(\"16\u00d738=608\", \"40\u00d732=1280\"),\n    @(\"18\u00d790=1620\", \"25\u00d721=525\"),\n    @(\"25\u00d716=400\", \"43\u00d734=1462\"),\n    @(\"85\u00d740=3400\", \"30\u00d724=720\"),\n    @(\"37\u00d711=407\", \"85\u00d735=2975\"),\n    @(\"34\u00d746=1564\", \"15\u00d711=165\"),\n    @(\"35\u00d747=1645\", \"94\u00d776=7144\"),\n    @(\"42\u00d797=4074\", \"64\u00d736=2304\"),\n    @(\"11\u00d734=374\", \"93\u00d771=6603\"),\n    @(\"62\u00d774=4588\", \"27\u00d794=2538\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
